# Applies the odds/stat updates for 2025-02-05 FlashScore weekly games workbook.
# Rows 3, 6, 7 and 8 in Sheet1 receive refreshed odds values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Zalaegerszegi vs Ferencvaros)
$ws.Range("G3").Value = 5.25
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 1.65
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 2.25
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.36
$ws.Range("U3").Value = 1.4
$ws.Range("V3").Value = 2.75
$ws.Range("Z3").Value = 26

# Row 6 (Celtic vs Dundee FC)
$ws.Range("G6").Value = 1.08
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 26
$ws.Range("K6").Value = 3.75
$ws.Range("L6").Value = 15
$ws.Range("O6").Value = 1.07
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 2.3
$ws.Range("U6").Value = 1.14
$ws.Range("V6").Value = 5.5
$ws.Range("W6").Value = 2.05
$ws.Range("X6").Value = 1.7
$ws.Range("Z6").Value = 8
$ws.Range("AB6").Value = 7
$ws.Range("AG6").Value = 34
$ws.Range("AH6").Value = 81
$ws.Range("AM6").Value = 151
$ws.Range("AN6").Value = 101
$ws.Range("AO6").Value = 351

# Row 7 (St. Gallen vs Lugano)
$ws.Range("M7").Value = 1.03
$ws.Range("O7").Value = 1.18
$ws.Range("T7").Value = 1.5

# Row 8 (Young Boys vs Yverdon)
$ws.Range("H8").Value = 5.25
$ws.Range("K8").Value = 2.6
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 17
$ws.Range("O8").Value = 1.17
$ws.Range("T8").Value = 1.57
$ws.Range("Z8").Value = 7.5
$ws.Range("AA8").Value = 8.5
$ws.Range("AP8").Value = 1.88
$ws.Range("AQ8").Value = 1.98
